$d = $word.ActiveDocument

# Locate the two paragraphs that need to be restructured:
#  - the paragraph ending in "...independente do banco de dados." which
#    currently carries the stray "_GoBack" bookmark
#  - the paragraph "e)DBAf)Usuário final leigog)..." which needs to be
#    split into "e)DBA" (its own paragraph, now carrying the bookmark)
#    and "f)Usuário final leigog)..." (separated by a blank paragraph,
#    matching the a)/b)/c)/d) items above it).
$pPermite = $null
$pDBA = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*independente do banco de dados*") {
        $pPermite = $p
    }
    if ($t -like "*DBAf*") {
        $pDBA = $p
        break
    }
}

$rng = $d.Range($pPermite.Range.Start, $pDBA.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="009960C2" w:rsidRDefault="009960C2"><w:r><w:t>Permite que a aplicação só use os dados quando necessário, ficando o código da aplicação independente do banco de dados.</w:t></w:r></w:p>
<w:p/>
<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>e)DBA</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="gramEnd"/></w:p>
<w:p w:rsidR="009960C2" w:rsidRDefault="009960C2"/>
<w:p w:rsidR="00701B22" w:rsidRDefault="00701B22"><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00701B22"><w:t>f)Usuário</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00701B22"><w:t xml:space="preserve"> final </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00701B22"><w:t>leigog</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00701B22"><w:t>)</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00701B22"><w:t>Transaçãoh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00701B22"><w:t xml:space="preserve">)Gerenciamento de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00701B22"><w:t>recuperaçãoi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00701B22"><w:t xml:space="preserve">)Gerenciamento de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00701B22"><w:t>concorrênciaj</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00701B22"><w:t>)Linguagem DDL.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng.InsertXML($xml)
